# O_Venovanie.docx dedication text update.
#
# 1. Extend the dedication line ("Eliške") with the rest of the
#    recipients, matching it to the other run's 12pt (sz/szCs=24) font
#    size.
# 2. The Word-managed "_GoBack" bookmark (last edit location) moves from
#    the dedication paragraph to the now-edited empty paragraph above it
#    — re-adding a bookmark with the same name relocates it, exactly as
#    Word itself does after an edit.

$d = $word.ActiveDocument

# --- 1. Append the rest of the dedication text -----------------------
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$para_range = $lastPara.Range
$para_range.InsertAfter(", Adíkovi, Táničke a Aničke.")

# Match the existing run's font size (12pt / half-points 24) across the
# whole paragraph, including the new run and the paragraph mark.
$para_range.Font.Size = 12
$para_range.Font.SizeBi = 12

# --- 2. Relocate the "_GoBack" bookmark to the edited paragraph ------
# (The 14th paragraph — an empty "NoSpacing" paragraph — is where the
# edit session left off; Word tracks this automatically with _GoBack.)
$target = $d.Paragraphs(14)
$d.Bookmarks.Add("_GoBack", $target.Range)
